# plantilla_tipo_explotacion.xlsx - reorganizacion completa
# - renombra la hoja "Datos" -> "tipo_explotacion"
# - renombra las cabeceras a minusculas sin acentos
# - quita el formato de cabecera (negrita blanca sobre relleno azul, centrado)
# - quita el ancho de columna personalizado (vuelve al ancho estandar)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "tipo_explotacion"

# New lower-case, unaccented header labels
$ws.Range("A1").Value = "codigo"
$ws.Range("B1").Value = "descripcion"
$ws.Range("C1").Value = "categoria"
$ws.Range("D1").Value = "comentario"

# Drop the bold/white-on-blue, centered header style entirely
$ws.Range("A1:D1").ClearFormats()

# Drop the custom 20-character column widths, back to the workbook default
$ws.Columns("A:D").ColumnWidth = $ws.StandardWidth
